# Google mobility and MYS data update
# - Extend the shared "Active cases exclude ICU" formula down through row 231
# - Normalise the style of the B:D,F:G cells on rows 225-227 (style 21 -> style 20)
# - Append 4 new days of MYS case data (rows 228-231), matching the formatting
#   used by the other populated data rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-style rows 225-227 data columns to match the rest of the table (style 20) ---
$ws.Range("B224:D224").Copy()
$ws.Range("B225:D227").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("F224:G224").Copy()
$ws.Range("F225:G227").PasteSpecial(-4122)   # xlPasteFormats

# --- 2. Prime the formatting for the new rows 228-231 by cloning row 225's formats ---
$ws.Range("A225:G225").Copy()
$ws.Range("A228:G228").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A229:G229").PasteSpecial(-4122)
$ws.Range("A230:G230").PasteSpecial(-4122)
$ws.Range("A231:G231").PasteSpecial(-4122)

# --- 3. New MYS case data ---
# Row 228: 2020-10-01 (serial 44105)
$ws.Range("A228").Value = 44105
$ws.Range("B228").Value = 260
$ws.Range("C228").Value = 1
$ws.Range("D228").Value = 1334
$ws.Range("E228").Formula = "=D228-F228"
$ws.Range("F228").Value = 20
$ws.Range("G228").Value = 3

# Row 229: 2020-10-02 (serial 44106)
$ws.Range("A229").Value = 44106
$ws.Range("B229").Value = 287
$ws.Range("C229").Value = 0
$ws.Range("D229").Value = 1540
$ws.Range("E229").Formula = "=D229-F229"
$ws.Range("F229").Value = 22
$ws.Range("G229").Value = 4

# Row 230: 2020-10-03 (serial 44107)
$ws.Range("A230").Value = 44107
$ws.Range("B230").Value = 317
$ws.Range("C230").Value = 0
$ws.Range("D230").Value = 1735
$ws.Range("E230").Formula = "=D230-F230"
$ws.Range("F230").Value = 29
$ws.Range("G230").Value = 4

# Row 231: 2020-10-04 (serial 44108)
$ws.Range("A231").Value = 44108
$ws.Range("B231").Value = 293
$ws.Range("C231").Value = 1
$ws.Range("D231").Value = 1961
$ws.Range("E231").Formula = "=D231-F231"
$ws.Range("F231").Value = 28
$ws.Range("G231").Value = 4
